# Update the "Metadata" sheet (Property / Value table).
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/family-id"
# Version
$meta.Range("B3").Value = "8.0.0"
# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# Update the "Elements" sheet.
$elements = $wb.Worksheets.Item("Elements")

# Row 2 = "Extension" (root element): the ele-1/ext-1 constraint text moved off
# of this row onto the new Extension.extension row, so clear it here.
$elements.Range("AI2").Value = ""

# Row 5 = "Extension.url": its Fixed Value mirrors the StructureDefinition URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/family-id"
